$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 830, pushing existing rows 830:894 down to 831:895
$ws.Rows.Item(830).Insert()

# Populate the new row 830 with the new data record
$ws.Cells.Item(830, 1).Value = 3
$ws.Cells.Item(830, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(830, 3).Value = "Coquimbo"
$ws.Cells.Item(830, 4).Value = 45265
$ws.Cells.Item(830, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(830, 5).Value = 5
$ws.Cells.Item(830, 6).Value = 100112032
$ws.Cells.Item(830, 7).Value = "Zapallo italiano"
$ws.Cells.Item(830, 8).Value = "Sin especificar"
$ws.Cells.Item(830, 9).Value = "Primera"
$ws.Cells.Item(830, 10).Value = 110
$ws.Cells.Item(830, 11).Value = 9000
$ws.Cells.Item(830, 12).Value = 9500
$ws.Cells.Item(830, 13).Value = 9273
$ws.Cells.Item(830, 14).Value = "$/caja 36 unidades"
$ws.Cells.Item(830, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(830, 16).Value = 258
$ws.Cells.Item(830, 17).Value = 36
$ws.Cells.Item(830, 18).Value = "Hortaliza"
